$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'243.44"
$ws.Range("G2").Value = "'12"
$ws.Range("D3").Value = "'25.03"
$ws.Range("G3").Value = "'12"
$ws.Range("D4").Value = "'5.141"
$ws.Range("G4").Value = "'12"
$ws.Range("D5").Value = "'0.05720"
$ws.Range("G5").Value = "'12"
$ws.Range("D6").Value = "'6.488"
$ws.Range("G6").Value = "'12"
$ws.Range("D7").Value = "'3.088"
$ws.Range("G7").Value = "'12"
$ws.Range("G8").Value = "'12"
$ws.Range("D9").Value = "'0.8430"
$ws.Range("G9").Value = "'12"
$ws.Range("G10").Value = "'12"
$ws.Range("D11").Value = "'0.06957"
$ws.Range("G11").Value = "'12"
$ws.Range("D12").Value = "'0.02839"
$ws.Range("G12").Value = "'12"
$ws.Range("D13").Value = "'0.09361"
$ws.Range("G13").Value = "'12"
$ws.Range("D14").Value = "'0.001527"
$ws.Range("G14").Value = "'12"
$ws.Range("D15").Value = "'0.0006001"
$ws.Range("E15").Value = "14OneONEWorstin24h"
$ws.Range("G15").Value = "'12"
$ws.Range("D16").Value = "'0.006163"
$ws.Range("G16").Value = "'12"
$ws.Range("D17").Value = "'3.501"
$ws.Range("G17").Value = "'12"
$ws.Range("D18").Value = "'2.091"
$ws.Range("G18").Value = "'12"
$ws.Range("D19").Value = "'0.3194"
$ws.Range("G19").Value = "'12"
$ws.Range("D20").Value = "'0.03117"
$ws.Range("G20").Value = "'12"
$ws.Range("D21").Value = "'0.1319"
$ws.Range("G21").Value = "'12"
$ws.Range("D22").Value = "'3.742"
$ws.Range("G22").Value = "'12"
$ws.Range("G23").Value = "'12"
$ws.Range("G24").Value = "'12"
$ws.Range("D25").Value = "'0.001237"
$ws.Range("G25").Value = "'12"
$ws.Range("D26").Value = "'0.004258"
$ws.Range("G26").Value = "'12"
$ws.Range("D27").Value = "'0.00009703"
$ws.Range("G27").Value = "'12"
$ws.Range("G28").Value = "'12"
$ws.Range("G29").Value = "'12"
$ws.Range("G30").Value = "'12"
$ws.Range("G31").Value = "'12"
$ws.Range("G32").Value = "'12"
$ws.Range("G33").Value = "'12"
$ws.Range("G34").Value = "'12"
$ws.Range("G35").Value = "'12"
$ws.Range("G36").Value = "'12"
$ws.Range("G37").Value = "'12"
$ws.Range("G38").Value = "'12"
$ws.Range("G39").Value = "'12"
$ws.Range("D40").Value = "'0.03622"
$ws.Range("G40").Value = "'12"
$ws.Range("D41").Value = "'0.006292"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("G41").Value = "'12"
$ws.Range("D42").Value = "'0.1049"
$ws.Range("G42").Value = "'12"
$ws.Range("D43").Value = "'0.003201"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("G43").Value = "'12"
$ws.Range("D44").Value = "'0.007381"
$ws.Range("G44").Value = "'12"
$ws.Range("D45").Value = "'0.00005291"
$ws.Range("G45").Value = "'12"
$ws.Range("G46").Value = "'12"
$ws.Range("D47").Value = "'0.1650"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("G47").Value = "'12"
$ws.Range("D48").Value = "'0.002300"
$ws.Range("G48").Value = "'12"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("G49").Value = "'12"
$ws.Range("D50").Value = "'0.0002001"
$ws.Range("G50").Value = "'12"
$ws.Range("G51").Value = "'12"
